$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 42611.88689814815

$ws.Range("B7").Value = 18
$ws.Range("C7").Value = 57
$ws.Range("D7").Value = 38
$ws.Range("E7").Value = 68
$ws.Range("F7").Value = 31
$ws.Range("G7").Value = 11829
$ws.Range("H7").Value = 23721
$ws.Range("I7").Value = 2660
$ws.Range("J7").Value = 337
$ws.Range("K7").Value = 227
$ws.Range("L7").Value = 15
$ws.Range("M7").Value = 7
$ws.Range("N7").Value = "Bag"
